$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 708.875
$ws.Range("I42").Value = 658.4
$ws.Range("J42").Value = 793
$ws.Range("K42").Value = 1975.2
$ws.Range("L42").Value = 2379
$ws.Range("M42").Value = -1745.2
$ws.Range("N42").Value = -2839

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2987.125
$ws.Range("I106").Value = 2987.125
$ws.Range("K106").Value = 2987.125
$ws.Range("M106").Value = -2356.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 33469748
$ws.Range("I132").Value = 41835830
$ws.Range("K132").Value = 125507490
$ws.Range("M132").Value = -125504960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2222.06
$ws.Range("I138").Value = 1094.4054
$ws.Range("J138").Value = 2884.3333
$ws.Range("K138").Value = 3283.2162
$ws.Range("L138").Value = 8652.999899999999
$ws.Range("M138").Value = 1856.7838
$ws.Range("N138").Value = -18932.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10545.84
$ws.Range("I32").Value = 6946.321
$ws.Range("J32").Value = 17359.215
$ws.Range("K32").Value = 6946.321
$ws.Range("L32").Value = 17359.215
$ws.Range("M32").Value = -6659.321
$ws.Range("N32").Value = -17933.215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6929416
$ws.Range("J63").Value = 4654.769
$ws.Range("L63").Value = 4654.769
$ws.Range("N63").Value = -6026.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6929416
$ws.Range("J66").Value = 4654.769
$ws.Range("L66").Value = 23273.845
$ws.Range("N66").Value = -30137.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 999.8570999999999
$ws.Range("J110").Value = 999.75
$ws.Range("L110").Value = 999.75
$ws.Range("N110").Value = -5089.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 43244.953
$ws.Range("J139").Value = 43244.953
$ws.Range("L139").Value = 43244.953
$ws.Range("N139").Value = -53524.953

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 29846.666
$ws.Range("J57").Value = 29846.666
$ws.Range("L57").Value = 29846.666
$ws.Range("N57").Value = -31286.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 23129.8
$ws.Range("I97").Value = 2454.7144
$ws.Range("K97").Value = 2454.7144
$ws.Range("M97").Value = -1463.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4457.3335
$ws.Range("I99").Value = 1563.3334
$ws.Range("K99").Value = 1563.3334
$ws.Range("M99").Value = -65.33339999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 29846.666
$ws.Range("J136").Value = 29846.666
$ws.Range("L136").Value = 29846.666
$ws.Range("N136").Value = -40046.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3302.5
$ws.Range("I31").Value = 1248.7037
$ws.Range("J31").Value = 6221.0527
$ws.Range("K31").Value = 1248.7037
$ws.Range("L31").Value = 6221.0527
$ws.Range("M31").Value = -953.7037
$ws.Range("N31").Value = -6811.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3302.5
$ws.Range("I34").Value = 1248.7037
$ws.Range("J34").Value = 6221.0527
$ws.Range("K34").Value = 1248.7037
$ws.Range("L34").Value = 6221.0527
$ws.Range("M34").Value = -1046.7037
$ws.Range("N34").Value = -6625.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5599.6
$ws.Range("I99").Value = 999.3333
$ws.Range("K99").Value = 999.3333
$ws.Range("M99").Value = 498.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2662.5908
$ws.Range("I122").Value = 2081.2144
$ws.Range("K122").Value = 6243.6432
$ws.Range("M122").Value = -3793.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5599.6
$ws.Range("I126").Value = 999.3333
$ws.Range("K126").Value = 2997.9999
$ws.Range("M126").Value = -527.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4474.864
$ws.Range("I132").Value = 3810.5625
$ws.Range("K132").Value = 11431.6875
$ws.Range("M132").Value = -8901.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5941.2964
$ws.Range("J134").Value = 3742.5715
$ws.Range("L134").Value = 11227.7145
$ws.Range("N134").Value = -16297.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49.47826
$ws.Range("I12").Value = 107.5
$ws.Range("K12").Value = 322.5
$ws.Range("M12").Value = -149.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6986.4707
$ws.Range("I68").Value = 867.2
$ws.Range("J68").Value = 15728.286
$ws.Range("K68").Value = 2601.6
$ws.Range("L68").Value = 47184.858
$ws.Range("M68").Value = -1790.6
$ws.Range("N68").Value = -48806.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 6986.4707
$ws.Range("I71").Value = 867.2
$ws.Range("J71").Value = 15728.286
$ws.Range("K71").Value = 7804.8
$ws.Range("L71").Value = 141554.574
$ws.Range("M71").Value = -3748.8
$ws.Range("N71").Value = -149666.574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1678.5834
$ws.Range("I81").Value = 562.6
$ws.Range("J81").Value = 2475.7144
$ws.Range("K81").Value = 1687.8
$ws.Range("L81").Value = 7427.1432
$ws.Range("M81").Value = -564.8000000000002
$ws.Range("N81").Value = -9673.143199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1678.5834
$ws.Range("I84").Value = 562.6
$ws.Range("J84").Value = 2475.7144
$ws.Range("K84").Value = 5063.400000000001
$ws.Range("L84").Value = 22281.4296
$ws.Range("M84").Value = 552.5999999999995
$ws.Range("N84").Value = -33513.4296

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 792.45715
$ws.Range("I113").Value = 700.63635
$ws.Range("J113").Value = 947.8461
$ws.Range("K113").Value = 2101.90905
$ws.Range("L113").Value = 2843.5383
$ws.Range("M113").Value = 68.09094999999979
$ws.Range("N113").Value = -7183.5383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4425.6
$ws.Range("I114").Value = 128
$ws.Range("J114").Value = 5500
$ws.Range("K114").Value = 384
$ws.Range("L114").Value = 16500
$ws.Range("M114").Value = 2870
$ws.Range("N114").Value = -23008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10000977
$ws.Range("J131").Value = 1002.9286
$ws.Range("L131").Value = 3008.7858
$ws.Range("N131").Value = -13088.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 32247.5
$ws.Range("J52").Value = 32247.5
$ws.Range("L52").Value = 32247.5
$ws.Range("N52").Value = -32765.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6025.355
$ws.Range("I70").Value = 5353.4546
$ws.Range("J70").Value = 7667.778
$ws.Range("K70").Value = 5353.4546
$ws.Range("L70").Value = 7667.778
$ws.Range("M70").Value = -5083.4546
$ws.Range("N70").Value = -8207.778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6025.355
$ws.Range("I73").Value = 5353.4546
$ws.Range("J73").Value = 7667.778
$ws.Range("K73").Value = 5353.4546
$ws.Range("L73").Value = 7667.778
$ws.Range("M73").Value = -4417.4546
$ws.Range("N73").Value = -9539.778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 28556.666
$ws.Range("J112").Value = 28556.666
$ws.Range("L112").Value = 28556.666
$ws.Range("N112").Value = -30772.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3488.79
$ws.Range("I126").Value = 2704.5762
$ws.Range("J126").Value = 4617.2925
$ws.Range("K126").Value = 8113.7286
$ws.Range("L126").Value = 13851.8775
$ws.Range("M126").Value = -5643.7286
$ws.Range("N126").Value = -18791.8775

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8930437
$ws.Range("I22").Value = 13890162
$ws.Range("J22").Value = 2930.1
$ws.Range("K22").Value = 13890162
$ws.Range("L22").Value = 2930.1
$ws.Range("M22").Value = -13889867
$ws.Range("N22").Value = -3520.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8930437
$ws.Range("I27").Value = 13890162
$ws.Range("J27").Value = 2930.1
$ws.Range("K27").Value = 13890162
$ws.Range("L27").Value = 2930.1
$ws.Range("M27").Value = -13890055
$ws.Range("N27").Value = -3144.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5810.242
$ws.Range("I40").Value = 5550.32
$ws.Range("J40").Value = 6622.5
$ws.Range("K40").Value = 5550.32
$ws.Range("L40").Value = 6622.5
$ws.Range("M40").Value = -5414.32
$ws.Range("N40").Value = -6894.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2280.9524
$ws.Range("I46").Value = 3533.6667
$ws.Range("J46").Value = 2072.1667
$ws.Range("K46").Value = 3533.6667
$ws.Range("L46").Value = 2072.1667
$ws.Range("M46").Value = -3345.6667
$ws.Range("N46").Value = -2448.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 281.53845
$ws.Range("I55").Value = 225.5
$ws.Range("K55").Value = 225.5
$ws.Range("M55").Value = -52.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 864.561
$ws.Range("I68").Value = 729.9231
$ws.Range("J68").Value = 3490
$ws.Range("K68").Value = 729.9231
$ws.Range("L68").Value = 3490
$ws.Range("M68").Value = 19.07690000000002
$ws.Range("N68").Value = -4988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 864.561
$ws.Range("I71").Value = 729.9231
$ws.Range("J71").Value = 3490
$ws.Range("K71").Value = 3649.6155
$ws.Range("L71").Value = 17450
$ws.Range("M71").Value = 94.38450000000012
$ws.Range("N71").Value = -24938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5407.615
$ws.Range("I122").Value = 4064.1428
$ws.Range("K122").Value = 12192.4284
$ws.Range("M122").Value = -9742.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3756.9
$ws.Range("I136").Value = 1926.8
$ws.Range("J136").Value = 5587
$ws.Range("K136").Value = 5780.4
$ws.Range("L136").Value = 16761
$ws.Range("M136").Value = -3230.4
$ws.Range("N136").Value = -21861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 47683.11
$ws.Range("J139").Value = 48643.5
$ws.Range("L139").Value = 48643.5
$ws.Range("N139").Value = -58923.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6800.4375
$ws.Range("I122").Value = 5317.25
$ws.Range("J122").Value = 11250
$ws.Range("K122").Value = 15951.75
$ws.Range("L122").Value = 33750
$ws.Range("M122").Value = -13501.75
$ws.Range("N122").Value = -38650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 569896.8
$ws.Range("I126").Value = 4190.5
$ws.Range("K126").Value = 12571.5
$ws.Range("M126").Value = -10101.5
